$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "description" column at C, shifting the existing "reason"
#     column (currently C) to D. Copy C1's header formatting onto D1 first,
#     then overwrite the header text in both columns. ---
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

$ws.Range("D1").Value = "reason"
$ws.Range("C1").Value = "description"

# --- Row 2: job id 3 (NLP Engineer) ---
$ws.Range("B2").Value = 94
$ws.Range("C2").Value = 'NLP Engineer (Remote): Collect and preprocess text corpora for language model training. Analyze data, develop and improve models. Skills: NLP, Pytorch, Computer Vision, Python.'
$ws.Range("D2").Value = 'The job of an NLP Engineer requires skills in NLP, Pytorch, Computer Vision, and Python. The candidate has experience in developing a pipeline for bias/toxicity detection in language models, utilizing GPT-2 and BERT, and implementing models for text analysis. The candidate''s skills align well with the requirements of the job, making it suitable for the candidate.'

# --- Row 3: job id changes 2 -> 4 (Flutter mobile app) ---
# Force the id to be stored as text (matches the rest of column A) without
# leaving a lingering cell style behind: write it as a formula producing the
# text "4", then convert the cell in place to a plain text value.
$ws.Range("A3").Formula = '="4"'
$ws.Range("A3").Copy()
$ws.Range("A3").PasteSpecial(-4163)   # xlPasteValues
$ws.Application.CutCopyMode = $false
$ws.Range("B3").Value = 41.5
$ws.Range("C3").Value = 'Application Development: Front-end coding for a sweat amino acid analysis app under a Scrum Master. Skills: Flutter, Dart, Android/iOS, Firestore, Firebase Authentication, Cloud Storage/Messaging, Mobile app architecture/design, Git.'
$ws.Range("D3").Value = 'The job involves front-end coding for a mobile app using Flutter and Dart, which aligns with your experience in ReactJS and JavaScript. The job also requires familiarity with mobile app architecture/design, which could be transferable skills from your previous projects. However, the job has a moderate score, suggesting that there may be other candidates with more relevant experience in Flutter and mobile app development.'

# --- Row 4: job id 1 (SDE Intern) ---
$ws.Range("B4").Value = 75
$ws.Range("C4").Value = 'SDE Intern: Remote MERN Stack Developer Internship, responsible for designing and developing web/mobile applications using MongoDB, ExpressJS, ReactJS, and NodeJS. Additional tasks include code maintenance, scalability, feature development, and product enhancement suggestions. Skills required: MongoDB, ReactJS, JavaScript, Web Development, NodeJS.'
$ws.Range("D4").Value = 'The job of a SDE Intern as a remote MERN Stack Developer requires skills in MongoDB, ReactJS, JavaScript, Web Development, and NodeJS. The candidate''s experience in developing web applications using ReactJS, JavaScript, and their proficiency in MongoDB and NodeJS make them suitable for this role. The high score of 75.0 indicates a significant overlap between the job requirements and the candidate''s skills.'

# --- Row 5 (new): job id 2 (Frontend Engineer Intern) ---
$ws.Range("A5").Formula = '="2"'
$ws.Range("A5").Copy()
$ws.Range("A5").PasteSpecial(-4163)   # xlPasteValues
$ws.Application.CutCopyMode = $false
$ws.Range("B5").Value = 71
$ws.Range("C5").Value = 'Frontend Engineer Intern - Work in a team to ensure consistent web design and user experience, optimize web pages, and maintain brand consistency. Requires excellent communication skills and proficiency in ReactJS, JavaScript, CSS, and NextJS. 3-month evaluative unpaid internship with potential return offers.'
$ws.Range("D5").Value = 'The job as a Frontend Engineer Intern involves working in a team to maintain brand consistency and optimize web pages. Your experience in developing user-friendly UI using ReactJS and JavaScript makes you suitable. However, the unpaid nature of the internship and the evaluative period might be slight drawbacks. Overall, it appears to be a moderate fit for your skills and experience.'

